$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 681.03
$ws.Range("I15").Value = 681.03
$ws.Range("K15").Value = 2043.09
$ws.Range("M15").Value = -1874.09
# row 132
$ws.Range("H132").Value = 2552.6826
$ws.Range("I132").Value = 2432.5178
$ws.Range("K132").Value = 7297.553400000001
$ws.Range("M132").Value = -4767.553400000001
# row 137
$ws.Range("H137").Value = 1178.6267
$ws.Range("I137").Value = 1117.921
$ws.Range("J137").Value = 1240.973
$ws.Range("K137").Value = 3353.763
$ws.Range("L137").Value = 3722.919
$ws.Range("M137").Value = -803.7629999999999
$ws.Range("N137").Value = -8822.919
# row 138
$ws.Range("H138").Value = 1090.09
$ws.Range("I138").Value = 564.1
$ws.Range("J138").Value = 1616.08
$ws.Range("K138").Value = 1692.3
$ws.Range("L138").Value = 4848.24
$ws.Range("M138").Value = 3447.7
$ws.Range("N138").Value = -15128.24
# row 141
$ws.Range("H141").Value = 2613.4043
$ws.Range("I141").Value = 985.4103
$ws.Range("J141").Value = 10549.875
$ws.Range("K141").Value = 2956.2309
$ws.Range("L141").Value = 31649.625
$ws.Range("M141").Value = 2223.7691
$ws.Range("N141").Value = -42009.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 3509.889
$ws.Range("I45").Value = 3999.8
$ws.Range("J45").Value = 2897.5
$ws.Range("K45").Value = 3999.8
$ws.Range("L45").Value = 2897.5
$ws.Range("M45").Value = -3622.8
$ws.Range("N45").Value = -3651.5
# row 61
$ws.Range("H61").Value = 10419021
$ws.Range("I61").Value = 12347729
$ws.Range("J61").Value = 3999.6
$ws.Range("K61").Value = 12347729
$ws.Range("L61").Value = 3999.6
$ws.Range("M61").Value = -12347517
$ws.Range("N61").Value = -4423.6
# row 74
$ws.Range("H74").Value = 1485.8064
$ws.Range("I74").Value = 742.3514
$ws.Range("J74").Value = 2586.12
$ws.Range("K74").Value = 742.3514
$ws.Range("L74").Value = 2586.12
$ws.Range("M74").Value = 131.6486
$ws.Range("N74").Value = -4334.12
# row 77
$ws.Range("H77").Value = 1485.8064
$ws.Range("I77").Value = 742.3514
$ws.Range("J77").Value = 2586.12
$ws.Range("K77").Value = 3711.757
$ws.Range("L77").Value = 12930.6
$ws.Range("M77").Value = 656.2429999999999
$ws.Range("N77").Value = -21666.6
# row 122
$ws.Range("H122").Value = 202149.4
$ws.Range("I122").Value = 251708.25
$ws.Range("J122").Value = 3914
$ws.Range("K122").Value = 755124.75
$ws.Range("L122").Value = 11742
$ws.Range("M122").Value = -752674.75
$ws.Range("N122").Value = -16642
# row 132
$ws.Range("H132").Value = 3260.9302
$ws.Range("I132").Value = 3475.875
$ws.Range("J132").Value = 2989.4211
$ws.Range("K132").Value = 10427.625
$ws.Range("L132").Value = 8968.263300000001
$ws.Range("M132").Value = -7897.625
$ws.Range("N132").Value = -14028.2633
# row 136
$ws.Range("H136").Value = 10419021
$ws.Range("I136").Value = 12347729
$ws.Range("J136").Value = 3999.6
$ws.Range("K136").Value = 37043187
$ws.Range("L136").Value = 11998.8
$ws.Range("M136").Value = -37040637
$ws.Range("N136").Value = -17098.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 2860.75
$ws.Range("I134").Value = 2617.9167
$ws.Range("J134").Value = 3589.25
$ws.Range("K134").Value = 7853.750100000001
$ws.Range("L134").Value = 10767.75
$ws.Range("M134").Value = -5318.750100000001
$ws.Range("N134").Value = -15837.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 513.3333
$ws.Range("I22").Value = 476.66666
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 476.66666
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -126.66666
$ws.Range("N22").Value = -1250
# row 31
$ws.Range("H31").Value = 3946.7222
$ws.Range("I31").Value = 1187.1951
$ws.Range("J31").Value = 7596.4194
$ws.Range("K31").Value = 1187.1951
$ws.Range("L31").Value = 7596.4194
$ws.Range("M31").Value = -892.1950999999999
$ws.Range("N31").Value = -8186.4194
# row 34
$ws.Range("H34").Value = 3946.7222
$ws.Range("I34").Value = 1187.1951
$ws.Range("J34").Value = 7596.4194
$ws.Range("K34").Value = 1187.1951
$ws.Range("L34").Value = 7596.4194
$ws.Range("M34").Value = -985.1950999999999
$ws.Range("N34").Value = -8000.4194
# row 100
$ws.Range("H100").Value = 45780
$ws.Range("J100").Value = 45780
$ws.Range("L100").Value = 45780
$ws.Range("N100").Value = -47944
# row 123
$ws.Range("H123").Value = 28213.334
$ws.Range("J123").Value = 28213.334
$ws.Range("L123").Value = 28213.334
$ws.Range("N123").Value = -38013.334
# row 132
$ws.Range("H132").Value = 5378316.5
$ws.Range("I132").Value = 1811.7222
$ws.Range("J132").Value = 12822707
$ws.Range("K132").Value = 5435.1666
$ws.Range("L132").Value = 38468121
$ws.Range("M132").Value = -2905.1666
$ws.Range("N132").Value = -38473181

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 1196.641
$ws.Range("I5").Value = 359.52173
$ws.Range("J5").Value = 2400
$ws.Range("K5").Value = 1078.56519
$ws.Range("L5").Value = 7200
$ws.Range("M5").Value = -966.56519
$ws.Range("N5").Value = -7424
# row 122
$ws.Range("H122").Value = 2822.9778
$ws.Range("I122").Value = 392.60715
$ws.Range("K122").Value = 3533.46435
$ws.Range("M122").Value = -1083.46435
# row 131
$ws.Range("H131").Value = 3109.9443
$ws.Range("I131").Value = 395.57144
$ws.Range("J131").Value = 4059.975
$ws.Range("K131").Value = 1186.71432
$ws.Range("L131").Value = 12179.925
$ws.Range("M131").Value = 3853.28568
$ws.Range("N131").Value = -22259.925
# row 135
$ws.Range("H135").Value = 1196.641
$ws.Range("I135").Value = 359.52173
$ws.Range("J135").Value = 2400
$ws.Range("K135").Value = 3235.69557
$ws.Range("L135").Value = 21600
$ws.Range("M135").Value = -700.6955699999999
$ws.Range("N135").Value = -26670

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 1734.1
$ws.Range("I102").Value = 1620.75
$ws.Range("K102").Value = 1620.75
$ws.Range("M102").Value = 1.25
# row 113
$ws.Range("H113").Value = 60666.41
$ws.Range("I113").Value = 92865.63
$ws.Range("J113").Value = 1634.5
$ws.Range("K113").Value = 92865.63
$ws.Range("L113").Value = 1634.5
$ws.Range("M113").Value = -90695.63
$ws.Range("N113").Value = -5974.5
# row 122
$ws.Range("H122").Value = 1949.875
$ws.Range("I122").Value = 1974.75
$ws.Range("J122").Value = 1925
$ws.Range("K122").Value = 5924.25
$ws.Range("L122").Value = 5775
$ws.Range("M122").Value = -3474.25
$ws.Range("N122").Value = -10675
# row 132
$ws.Range("H132").Value = 3311.606
$ws.Range("I132").Value = 2829.6843
$ws.Range("J132").Value = 3965.6428
$ws.Range("K132").Value = 8489.052899999999
$ws.Range("L132").Value = 11896.9284
$ws.Range("M132").Value = -5959.052899999999
$ws.Range("N132").Value = -16956.9284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2169
$ws.Range("I7").Value = 2013.25
$ws.Range("J7").Value = 2565.4546
$ws.Range("K7").Value = 2013.25
$ws.Range("L7").Value = 2565.4546
$ws.Range("M7").Value = -1901.25
$ws.Range("N7").Value = -2789.4546
# row 40
$ws.Range("H40").Value = 44129.25
$ws.Range("J40").Value = 2868.1667
$ws.Range("L40").Value = 2868.1667
$ws.Range("N40").Value = -3140.1667
# row 61
$ws.Range("H61").Value = 2377.862
$ws.Range("I61").Value = 2062.7827
$ws.Range("J61").Value = 3585.6667
$ws.Range("K61").Value = 2062.7827
$ws.Range("L61").Value = 3585.6667
$ws.Range("M61").Value = -1860.7827
$ws.Range("N61").Value = -3989.6667
# row 113
$ws.Range("H113").Value = 2377.862
$ws.Range("I113").Value = 2062.7827
$ws.Range("J113").Value = 3585.6667
$ws.Range("K113").Value = 2062.7827
$ws.Range("L113").Value = 3585.6667
$ws.Range("M113").Value = 107.2172999999998
$ws.Range("N113").Value = -7925.6667
# row 122
$ws.Range("H122").Value = 4668.864
$ws.Range("I122").Value = 3716
$ws.Range("J122").Value = 4949.1177
$ws.Range("K122").Value = 11148
$ws.Range("L122").Value = 14847.3531
$ws.Range("M122").Value = -8698
$ws.Range("N122").Value = -19747.3531
# row 126
$ws.Range("H126").Value = 2169
$ws.Range("I126").Value = 2013.25
$ws.Range("J126").Value = 2565.4546
$ws.Range("K126").Value = 6039.75
$ws.Range("L126").Value = 7696.3638
$ws.Range("M126").Value = -3569.75
$ws.Range("N126").Value = -12636.3638
# row 136
$ws.Range("H136").Value = 3789377.8
$ws.Range("I136").Value = 1309.7858
$ws.Range("J136").Value = 10418497
$ws.Range("K136").Value = 3929.3574
$ws.Range("L136").Value = 31255491
$ws.Range("M136").Value = -1379.3574
$ws.Range("N136").Value = -31260591

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 113
$ws.Range("H113").Value = 1059.7241
$ws.Range("I113").Value = 1077.68
$ws.Range("K113").Value = 3233.04
$ws.Range("M113").Value = -1063.04
# row 122
$ws.Range("H122").Value = 3001.2163
$ws.Range("I122").Value = 2757.6
$ws.Range("J122").Value = 3508.75
$ws.Range("K122").Value = 8272.799999999999
$ws.Range("L122").Value = 10526.25
$ws.Range("M122").Value = -5822.799999999999
$ws.Range("N122").Value = -15426.25
# row 132
$ws.Range("H132").Value = 10060228
$ws.Range("I132").Value = 3978.077
$ws.Range("J132").Value = 18230930
$ws.Range("K132").Value = 11934.231
$ws.Range("L132").Value = 54692790
$ws.Range("M132").Value = -9404.231
$ws.Range("N132").Value = -54697850
# row 136
$ws.Range("H136").Value = 2087.3333
$ws.Range("I136").Value = 2179.5356
$ws.Range("J136").Value = 1920.7742
$ws.Range("K136").Value = 6538.6068
$ws.Range("L136").Value = 5762.3226
$ws.Range("M136").Value = -3988.6068
$ws.Range("N136").Value = -10862.3226

